$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of existing header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I and J columns for each data row (2-31)
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 4
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 6
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 7
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 6
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 6
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 8
$ws.Range("I11").Value = 8
$ws.Range("J11").Value = 8
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 5
$ws.Range("I13").Value = 6
$ws.Range("J13").Value = 6
$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 7
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 8
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 6
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 9
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 8
$ws.Range("I19").Value = 7
$ws.Range("J19").Value = 7
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 8
$ws.Range("I21").Value = 6
$ws.Range("J21").Value = 6
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 8
$ws.Range("I23").Value = 5
$ws.Range("J23").Value = 5
$ws.Range("I24").Value = 8
$ws.Range("J24").Value = 8
$ws.Range("I25").Value = 6
$ws.Range("J25").Value = 7
$ws.Range("I26").Value = 6
$ws.Range("J26").Value = 6
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 7
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 5
$ws.Range("I29").Value = 6
$ws.Range("J29").Value = 6
$ws.Range("I30").Value = 7
$ws.Range("J30").Value = 7
$ws.Range("I31").Value = 7
$ws.Range("J31").Value = 7
